# Committing fixed search test cases
# Replace placeholder Jira ids (TBD-N / TDB-8 / OPQA-610) in the "Jira id"
# column (B) of the "Test Cases" sheet with their real ticket numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("B26").Value = "OPQA-1434"
$ws.Range("B27").Value = "OPQA-1435"
$ws.Range("B28").Value = "OPQA-1436"
$ws.Range("B29").Value = "OPQA-1437"
$ws.Range("B30").Value = "OPQA-1438"
$ws.Range("B31").Value = "OPQA-1439"
$ws.Range("B32").Value = "OPQA-1440"
$ws.Range("B33").Value = "OPQA-1441"
$ws.Range("B34").Value = "OPQA-1442"
$ws.Range("B35").Value = "OPQA-1443"
$ws.Range("B36").Value = "OPQA-1444"
$ws.Range("B37").Value = "OPQA-1445"
$ws.Range("B38").Value = "OPQA-1447"
$ws.Range("B39").Value = "OPQA-1449"
$ws.Range("B40").Value = "OPQA-1450"
$ws.Range("B41").Value = "OPQA-1452"
$ws.Range("B42").Value = "OPQA-1453"
$ws.Range("B43").Value = "OPQA-1455"
$ws.Range("B44").Value = "OPQA-1456"
$ws.Range("B45").Value = "OPQA-1501"
